$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value looks like a plain number (single decimal point)
# must be kept as text, matching the workbook's original inline-string cells.
# Setting NumberFormat to "@" (Text) before assigning the value prevents Excel
# from silently converting the text into a numeric value.
$textCells = @(
    "D5",
    "D6",
    "D10",
    "D12",
    "D16",
    "D19",
    "D22",
    "D24",
    "D26",
    "D27",
    "D32",
    "D34",
    "D35",
    "D37",
    "D39",
    "D40",
    "D41",
    "D42",
    "D43",
    "D44",
    "D45",
    "D47",
    "D48",
    "D50",
    "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "59.222.34"
$ws.Range("E2").Value = "  +2.39%  "
$ws.Range("D3").Value = "2.596.46"
$ws.Range("E3").Value = "  +1.56%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "534.38"
$ws.Range("E5").Value = "  +3.81%  "
$ws.Range("D6").Value = "140.30"
$ws.Range("E6").Value = "  +1.79%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +1.68%  "
$ws.Range("D9").Value = "2.607.95"
$ws.Range("E9").Value = "  +1.77%  "
$ws.Range("D10").Value = "6.47"
$ws.Range("E10").Value = "  +0.21%  "
$ws.Range("E11").Value = "  +4.07%  "
$ws.Range("D12").Value = "0.335"
$ws.Range("E12").Value = "  +3.60%  "
$ws.Range("E13").Value = "  +2.69%  "
$ws.Range("D14").Value = "3.055.90"
$ws.Range("E14").Value = "  +1.71%  "
$ws.Range("D15").Value = "59.159.47"
$ws.Range("E15").Value = "  +2.32%  "
$ws.Range("D16").Value = "20.49"
$ws.Range("E16").Value = "  +2.83%  "
$ws.Range("D17").Value = "2.598.60"
$ws.Range("E17").Value = "  +1.16%  "
$ws.Range("E18").Value = "  +2.36%  "
$ws.Range("D19").Value = "345.60"
$ws.Range("E19").Value = "  +4.22%  "
$ws.Range("E20").Value = "  +1.92%  "
$ws.Range("E21").Value = "  +0.92%  "
$ws.Range("D22").Value = "6.34"
$ws.Range("E22").Value = "  +0.82%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").Value = "67.14"
$ws.Range("E24").Value = "  +2.71%  "
$ws.Range("E25").Value = "  +3.13%  "
$ws.Range("D26").Value = "0.406"
$ws.Range("E26").Value = "  +2.84%  "
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("E28").Value = "  +4.54%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "0.0₃0734"
$ws.Range("E30").Value = "  +4.65%  "
$ws.Range("E31").Value = "  +5.29%  "
$ws.Range("D32").Value = "5.80"
$ws.Range("E32").Value = "  -1.23%  "
$ws.Range("E33").Value = "  +2.00%  "
$ws.Range("D34").Value = "149.46"
$ws.Range("E34").Value = "  +0.45%  "
$ws.Range("D35").Value = "3.97"
$ws.Range("E35").Value = "  +3.07%  "
$ws.Range("E36").Value = "  +1.56%  "
$ws.Range("D37").Value = "36.96"
$ws.Range("E37").Value = "  +2.98%  "
$ws.Range("E38").Value = "  +5.38%  "
$ws.Range("D39").Value = "0.835"
$ws.Range("E39").Value = "  +1.37%  "
$ws.Range("D40").Value = "0.823"
$ws.Range("E40").Value = "  +0.70%  "
$ws.Range("D41").Value = "3.53"
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "0.998"
$ws.Range("E42").Value = "  +0.20%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").Value = "276.43"
$ws.Range("E43").Value = "  +3.37%  "
$ws.Range("D44").Value = "0.598"
$ws.Range("E44").Value = "  +3.11%  "
$ws.Range("D45").Value = "10.76"
$ws.Range("E45").Value = "  +0.50%  "
$ws.Range("E46").Value = "  +3.01%  "
$ws.Range("D47").Value = "0.0521"
$ws.Range("E47").Value = "  +2.27%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "0.0223"
$ws.Range("E48").Value = "  +3.74%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "1.938.69"
$ws.Range("E49").Value = "  -0.29%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "18.29"
$ws.Range("E50").Value = "  +4.59%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "4.51"
$ws.Range("E51").Value = "  +3.25%  "
